$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 554
$ws1.Range("F5").Value = 305
$ws1.Range("F7").Value = 254
$ws1.Range("F8").Value = 2344
$ws1.Range("F10").Value = 5897
$ws1.Range("F11").Value = 147

# Sheet "全部类型" (All types) - fourth sheet, contains combined data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 554
$ws4.Range("F6").Value = 305
$ws4.Range("F8").Value = 254
$ws4.Range("F11").Value = 2344
$ws4.Range("F13").Value = 5897
$ws4.Range("F14").Value = 147
